$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving its original Text cell type
# (the source data stores numbers as text) and without leaving behind a
# changed/quote-prefixed number format on the cell.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("C2") "1013"
Set-TextValue $ws.Range("D2") "2242686.79"
Set-TextValue $ws.Range("C4") "1055"
Set-TextValue $ws.Range("D4") "3804136.47"
Set-TextValue $ws.Range("C6") "697"
Set-TextValue $ws.Range("D6") "2319021.78"
Set-TextValue $ws.Range("C20") "185"
Set-TextValue $ws.Range("D20") "470859.00"
Set-TextValue $ws.Range("C21") "343"
Set-TextValue $ws.Range("D21") "1240274.00"
Set-TextValue $ws.Range("C22") "163"
Set-TextValue $ws.Range("D22") "488812.39"
Set-TextValue $ws.Range("C24") "12"
Set-TextValue $ws.Range("D24") "56250.00"
Set-TextValue $ws.Range("C29") "303"
Set-TextValue $ws.Range("D29") "780236.89"
Set-TextValue $ws.Range("C41") "213"
Set-TextValue $ws.Range("D41") "600155.23"
Set-TextValue $ws.Range("C42") "101"
Set-TextValue $ws.Range("D42") "484699.98"
Set-TextValue $ws.Range("C43") "144"
Set-TextValue $ws.Range("D43") "591472.25"
Set-TextValue $ws.Range("C44") "5"
Set-TextValue $ws.Range("D44") "19519.00"
Set-TextValue $ws.Range("C46") "425"
Set-TextValue $ws.Range("D46") "1175854.43"
Set-TextValue $ws.Range("C48") "660"
Set-TextValue $ws.Range("D48") "2729876.99"
Set-TextValue $ws.Range("C49") "453"
Set-TextValue $ws.Range("D49") "1641996.50"
Set-TextValue $ws.Range("C55") "4058"
Set-TextValue $ws.Range("D55") "14294025.38"
Set-TextValue $ws.Range("C77") "952"
Set-TextValue $ws.Range("D77") "3353064.26"
Set-TextValue $ws.Range("C78") "535"
Set-TextValue $ws.Range("D78") "1777407.38"
Set-TextValue $ws.Range("C93") "719"
Set-TextValue $ws.Range("D93") "1750599.94"
Set-TextValue $ws.Range("C96") "1155"
Set-TextValue $ws.Range("D96") "4005663.23"
Set-TextValue $ws.Range("C98") "1098"
Set-TextValue $ws.Range("D98") "3517411.30"
Set-TextValue $ws.Range("C99") "16"
Set-TextValue $ws.Range("D99") "46345.28"
Set-TextValue $ws.Range("C100") "54"
Set-TextValue $ws.Range("D100") "206835.52"
